# Manual dislocation uploading 2021/11/06 21:00
# The monthly "Sentinel" export window rolled forward by one month: the
# duplicated Sep/Oct block (rows 2-31 and 32-61, which both repeated the same
# 2021-09-01..2021-09-30 date stamps) is replaced with one continuous run of
# real daily records from 2021-10-01 through 2021-11-30 (rows 2-62), all
# shipped from the same route (Aktogay -> Dostyk, cargo "MED'").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fromStation = "Актогай"
$toStation   = "Достык (эксп.)"
$cargoName   = "МЕДЬ"

# Non-zero CarAmount overrides for specific rows; every other row in the
# 2..62 range settles to 0.
$carAmountOverrides = @{ 16 = 46; 31 = 46; 52 = 50 }

# Row 2 starts at serial date 44470 (2021-10-01) and the whole block is one
# unbroken run of consecutive days through row 62 (44530 / 2021-11-30).
$startSerial = 44470

for ($r = 2; $r -le 62; $r++) {
    $serial = $startSerial + ($r - 2)

    if ($carAmountOverrides.ContainsKey($r)) {
        $carAmount = $carAmountOverrides[$r]
    } else {
        $carAmount = 0
    }

    $ws.Cells.Item($r, 1).Value = $serial
    $ws.Cells.Item($r, 2).Value = $carAmount
    $ws.Cells.Item($r, 3).Value = $fromStation
    $ws.Cells.Item($r, 4).Value = $toStation
    $ws.Cells.Item($r, 5).Value = $cargoName
}

# View-state nudge: user had scrolled further down and selected B31 by the
# time this batch was saved.
[void]$ws.Range("B31").Select()

Write-Host "KAL.xlsx dislocation upload applied (rows 2-62 reset to 2021-10-01..2021-11-30)."
